# Generate Report for Handoff
#
# Replaces every occurrence of the old handoff-file UUID
# (2a3e7af6-5866-4363-b48b-53236c73f29f) with the new one
# (ba8290e7-589c-44af-9ee0-da02b12931cf), and refreshes the associated
# handoff/handback timestamps, on all three worksheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "2a3e7af6-5866-4363-b48b-53236c73f29f"
$newGuid = "ba8290e7-589c-44af-9ee0-da02b12931cf"

$oldZhXlf = "$oldGuid.afaf435430058c3b197c85fe13fc85515a48f83c.zh-cn.xlf"
$newZhXlf = "$newGuid.e0091d0906d7d380ebfa7bc3b0e53fb7c5a97800.zh-cn.xlf"

$oldDeXlf = "$oldGuid.afaf435430058c3b197c85fe13fc85515a48f83c.de-de.xlf"
$newDeXlf = "$newGuid.e0091d0906d7d380ebfa7bc3b0e53fb7c5a97800.de-de.xlf"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5900579d935581347ad9eb28ab517374bdd89e9/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-21 05:03:19"

# Recreate the hyperlink on B2 so its display text matches the new file
# name while keeping the original link target untouched.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-21 05:03:14"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-21 05:03:19"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
